$d = $word.ActiveDocument

$replacements = @(
    @("23×46=", "67×62="),
    @("44×85=", "78×63="),
    @("49×15=", "38×66="),
    @("48×48=", "63×82="),
    @("60×62=", "90×76="),
    @("41×24=", "77×38="),
    @("34×36=", "75×67="),
    @("43×82=", "89×13="),
    @("78×91=", "56×99="),
    @("23×52=", "64×75="),
    @("47×61=", "14×55="),
    @("19×83=", "73×17="),
    @("39×87=", "66×52="),
    @("58×75=", "12×64="),
    @("36×13=", "73×89="),
    @("28×85=", "96×73="),
    @("19×42=", "65×98="),
    @("37×66=", "19×50="),
    @("11×87=", "51×92="),
    @("80×86=", "24×87="),
    @("48×13=", "79×37="),
    @("68×14=", "78×48="),
    @("56×20=", "75×35="),
    @("68×60=", "29×61="),
    @("47×35=", "82×79=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
